$wb = $excel.ActiveWorkbook

# Sheet 1 ("展览") and Sheet 4 ("全部类型") both contain the same data table
# and both need the "想去人数" (want-to-go count) column F updated for
# rows 3, 4, and 6.
$sheetIndexes = @(1, 4)

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Range("F3").Value = 2140
    $ws.Range("F4").Value = 300
    $ws.Range("F6").Value = 6393
}
